$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 199, shifting the existing rows 199-202 down to 200-203.
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with the new weekly price-report entry.
$ws.Cells.Item(199, 1).Value  = 3
$ws.Cells.Item(199, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(199, 3).Value  = "Coquimbo"
$ws.Cells.Item(199, 4).Value  = 44448
$ws.Cells.Item(199, 5).Value  = 5
$ws.Cells.Item(199, 6).Value  = 100112031
$ws.Cells.Item(199, 7).Value  = "Poroto verde"
$ws.Cells.Item(199, 8).Value  = "Magnum"
$ws.Cells.Item(199, 9).Value  = "Primera"
$ws.Cells.Item(199, 10).Value = 85
$ws.Cells.Item(199, 11).Value = 33000
$ws.Cells.Item(199, 12).Value = 34000
$ws.Cells.Item(199, 13).Value = 33471
$ws.Cells.Item(199, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(199, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(199, 16).Value = 1339
$ws.Cells.Item(199, 17).Value = 25
$ws.Cells.Item(199, 18).Value = "Hortaliza"
